$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Row 11 / column B ("R40" rule-name cell) is changed from the text "R40"
# to the text "1". Build it as a text formula first and paste back as a
# value so the result keeps its original General-text formatting/style
# (an ordinary .Value assignment of a numeric-looking string would be
# auto-converted to a number by Excel's type inference).
$ws.Range("B11").Formula = '="1"'
$ws.Range("B11").Copy()
$ws.Range("B11").PasteSpecial(-4163)
